# Replaced buck converter inductor (SRR1210-221M / 220u -> MPL-SE4030-220 / 22u)
# Update the inductance input (E20) and inductor ESR input (A36) on each
# regulator sheet, and swap the part number / value / price in the parts
# list row that carries the inductor ("L" purpose column).

$wb = $excel.ActiveWorkbook

# --- "5V" sheet (inductor row = 78) ---
$ws1 = $wb.Worksheets.Item("5V")
$ws1.Range("E20").Value = 0.000022
$ws1.Range("A36").Value = 0.219
$ws1.Range("B78").Value = "22u"
$ws1.Range("A78").Value = "MPL-SE4030-220"
$ws1.Range("D78").Value = 0.56

# --- "3V7" sheet (inductor row = 80) ---
$ws2 = $wb.Worksheets.Item("3V7")
$ws2.Range("E20").Value = 0.000022
$ws2.Range("A36").Value = 0.219
$ws2.Range("B80").Value = "22u"
$ws2.Range("A80").Value = "MPL-SE4030-220"
$ws2.Range("D80").Value = 0.56

# --- "3V3" sheet (inductor row = 80) ---
$ws3 = $wb.Worksheets.Item("3V3")
$ws3.Range("E20").Value = 0.000022
$ws3.Range("A36").Value = 0.219
$ws3.Range("B80").Value = "22u"
$ws3.Range("A80").Value = "MPL-SE4030-220"
$ws3.Range("D80").Value = 0.56

# --- Restore the view/selection state seen in the saved workbook ---
$ws2.Activate()
$ws2.Range("A80:D80").Select()

$ws3.Activate()
$ws3.Range("A80:D80").Select()

$ws1.Activate()
$ws1.Range("C85").Select()
